# ---------------------------------------------------------------------------
# Applies the "added spanish and basque, other updates" edit to the house
# rules document:
#   1. The old "If you have any questions please contact the owners on: "
#      paragraph becomes empty, four more blank Subtitle paragraphs are
#      added, and a new final paragraph with the "Please inform Alan and
#      Jen..." text (and a lastRenderedPageBreak marker) is appended.
#   2. "If it's urgent contact our agent Oscar on" is split into "Or" +
#      " our agent Oscar on" (two runs), keeping the trailing ":" run.
#   3. A new empty paragraph is added right before the table.
#   4. The stray <w:lastRenderedPageBreak/> on the picture run is removed
#      (it moved to the new text run added in step 1).
#   5. The two runs making up "Empty the dehumidifier's water tank
#      regularly" + "." are merged into a single run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like ("*" + $needle + "*")) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Step 1: "If you have any questions..." paragraph -> 5 Subtitle paragraphs
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "If you have any questions"
$target = $d.Paragraphs($idx)

$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Please inform Alan and Jen, immediately if there are any problems at</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Step 2: "If it's urgent contact our agent Oscar on" paragraph -> "Or" + " our agent Oscar on" + ":"
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "our agent Oscar on"
$target = $d.Paragraphs($idx)

$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t>Or</w:t></w:r><w:r><w:t xml:space="preserve"> our agent Oscar on</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Step 3: insert an empty paragraph right before the table (after "+34 610 79 37 48")
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "+34 610 79 37 48"
$target = $d.Paragraphs($idx)
$target.Range.InsertParagraphAfter()

$newIdx = $idx + 1
$newPara = $d.Paragraphs($newIdx)
$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$newPara.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Step 4: remove <w:lastRenderedPageBreak/> from the picture run (first table
# cell, first paragraph).
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)
$picCell = $tbl.Cell(1, 1)
$picPara = $picCell.Range.Paragraphs(1)

$xml4 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0050F2B5" w14:textId="7E477D2F" w:rsidR="00D241E5" w:rsidRDefault="00D241E5" w:rsidP="009424E8"><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251669504" behindDoc="1" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="179A7EC5" wp14:editId="3225616B"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:posOffset>-1905</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>172720</wp:posOffset></wp:positionV><wp:extent cx="2929255" cy="1646555"/><wp:effectExtent l="0" t="0" r="4445" b="0"/><wp:wrapTight wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="0"/><wp:lineTo x="0" y="21242"/><wp:lineTo x="21492" y="21242"/><wp:lineTo x="21492" y="0"/><wp:lineTo x="0" y="0"/></wp:wrapPolygon></wp:wrapTight><wp:docPr id="298525180" name="Picture 1" descr="A close up of a microwave&#xA;&#xA;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="298525180" name="Picture 1" descr="A close up of a microwave&#xA;&#xA;Description automatically generated"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="2929255" cy="1646555"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$picPara.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# Step 5: merge the two runs of "Empty the dehumidifier's water tank
# regularly" + "." into a single run (inside the table).
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex $d "water tank regularly"
$target = $d.Paragraphs($idx)

$xml5 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Empty the dehumidifier&#8217;s water tank regularly.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml5)

Write-Host "Edit complete."
